$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "PESSOA 1"
$ws.Range("B2").Value = 99999999999
$ws.Range("B4").Select()
